$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '72.597.91'
$ws.Range('E2').Value = '  +0.39%  '

# Row 3
$ws.Range('D3').Value = '3.982.41'
$ws.Range('E3').Value = '  -0.94%  '

# Row 4
$ws.Range('E4').Value = '  -0.25%  '

# Row 5
$ws.Range('D5').Formula = '="587.12"'
$ws.Range('D5').Copy()
$ws.Range('D5').PasteSpecial(-4163)
$ws.Range('E5').Value = '  +11.41%  '

# Row 6
$ws.Range('D6').Formula = '="151.60"'
$ws.Range('D6').Copy()
$ws.Range('D6').PasteSpecial(-4163)
$ws.Range('E6').Value = '  +0.75%  '

# Row 7
$ws.Range('D7').Formula = '="0.679"'
$ws.Range('D7').Copy()
$ws.Range('D7').PasteSpecial(-4163)
$ws.Range('E7').Value = '  -3.13%  '

# Row 8
$ws.Range('D8').Formula = '="0.999"'
$ws.Range('D8').Copy()
$ws.Range('D8').PasteSpecial(-4163)
$ws.Range('E8').Value = '  -0.36%  '

# Row 9
$ws.Range('D9').Formula = '="0.751"'
$ws.Range('D9').Copy()
$ws.Range('D9').PasteSpecial(-4163)
$ws.Range('E9').Value = '  +0.30%  '

# Row 10
$ws.Range('D10').Formula = '="0.168"'
$ws.Range('D10').Copy()
$ws.Range('D10').PasteSpecial(-4163)
$ws.Range('E10').Value = '  -2.01%  '

# Row 11
$ws.Range('D11').Formula = '="53.64"'
$ws.Range('D11').Copy()
$ws.Range('D11').PasteSpecial(-4163)
$ws.Range('E11').Value = '  +6.94%  '

# Row 12
$ws.Range('D12').Formula = '="0.0000318"'
$ws.Range('D12').Copy()
$ws.Range('D12').PasteSpecial(-4163)
$ws.Range('E12').Value = '  -1.68%  '

# Row 13
$ws.Range('D13').Formula = '="10.86"'
$ws.Range('D13').Copy()
$ws.Range('D13').PasteSpecial(-4163)
$ws.Range('E13').Value = '  +1.33%  '

# Row 14
$ws.Range('D14').Value = '4.622.54'
$ws.Range('E14').Value = '  -1.15%  '

# Row 15
$ws.Range('D15').Value = '3.981.06'
$ws.Range('E15').Value = '  -1.50%  '

# Row 16
$ws.Range('D16').Formula = '="1.28"'
$ws.Range('D16').Copy()
$ws.Range('D16').PasteSpecial(-4163)
$ws.Range('E16').Value = '  +7.90%  '

# Row 17
$ws.Range('D17').Formula = '="14.03"'
$ws.Range('D17').Copy()
$ws.Range('D17').PasteSpecial(-4163)
$ws.Range('E17').Value = '  -0.06%  '

# Row 18
$ws.Range('D18').Formula = '="20.47"'
$ws.Range('D18').Copy()
$ws.Range('D18').PasteSpecial(-4163)
$ws.Range('E18').Value = '  -0.71%  '

# Row 19
$ws.Range('E19').Value = '  -0.55%  '

# Row 20
$ws.Range('D20').Value = '72.445.74'
$ws.Range('E20').Value = '  +0.33%  '

# Row 21
$ws.Range('D21').Formula = '="431.90"'
$ws.Range('D21').Copy()
$ws.Range('D21').PasteSpecial(-4163)
$ws.Range('E21').Value = '  +0.91%  '

# Row 22
$ws.Range('D22').Formula = '="4.71"'
$ws.Range('D22').Copy()
$ws.Range('D22').PasteSpecial(-4163)
$ws.Range('E22').Value = '  +12.89%  '

# Row 23
$ws.Range('D23').Formula = '="95.80"'
$ws.Range('D23').Copy()
$ws.Range('D23').PasteSpecial(-4163)
$ws.Range('E23').Value = '  -1.62%  '

# Row 24
$ws.Range('D24').Formula = '="3.47"'
$ws.Range('D24').Copy()
$ws.Range('D24').PasteSpecial(-4163)
$ws.Range('E24').Value = '  -0.44%  '

# Row 25
$ws.Range('B25').Value = 'Toncoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D25').Formula = '="4.47"'
$ws.Range('D25').Copy()
$ws.Range('D25').PasteSpecial(-4163)
$ws.Range('E25').Value = '  +20.82%  '

# Row 26
$ws.Range('B26').Value = 'InternetComputer(DFINITY)'
$ws.Range('C26').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D26').Formula = '="14.23"'
$ws.Range('D26').Copy()
$ws.Range('D26').PasteSpecial(-4163)
$ws.Range('E26').Value = '  -0.51%  '

# Row 27
$ws.Range('D27').Formula = '="11.41"'
$ws.Range('D27').Copy()
$ws.Range('D27').PasteSpecial(-4163)
$ws.Range('E27').Value = '  +2.55%  '

# Row 28
$ws.Range('D28').Formula = '="10.62"'
$ws.Range('D28').Copy()
$ws.Range('D28').PasteSpecial(-4163)
$ws.Range('E28').Value = '  -2.25%  '

# Row 29
$ws.Range('D29').Formula = '="5.92"'
$ws.Range('D29').Copy()
$ws.Range('D29').PasteSpecial(-4163)
$ws.Range('E29').Value = '  +0.93%  '

# Row 30
$ws.Range('D30').Formula = '="36.41"'
$ws.Range('D30').Copy()
$ws.Range('D30').PasteSpecial(-4163)
$ws.Range('E30').Value = '  -0.75%  '

# Row 31
$ws.Range('D31').Formula = '="7.99"'
$ws.Range('D31').Copy()
$ws.Range('D31').PasteSpecial(-4163)
$ws.Range('E31').Value = '  +8.19%  '

# Row 32
$ws.Range('D32').Formula = '="49.85"'
$ws.Range('D32').Copy()
$ws.Range('D32').PasteSpecial(-4163)
$ws.Range('E32').Value = '  +2.79%  '

# Row 33
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').Formula = '="13.51"'
$ws.Range('D33').Copy()
$ws.Range('D33').PasteSpecial(-4163)
$ws.Range('E33').Value = '  +0.77%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Formula = '="0.132"'
$ws.Range('D34').Copy()
$ws.Range('D34').PasteSpecial(-4163)
$ws.Range('E34').Value = '  +0.25%  '

# Row 35
$ws.Range('D35').Formula = '="681.31"'
$ws.Range('D35').Copy()
$ws.Range('D35').PasteSpecial(-4163)
$ws.Range('E35').Value = '  -0.14%  '

# Row 36
$ws.Range('D36').Formula = '="69.66"'
$ws.Range('D36').Copy()
$ws.Range('D36').PasteSpecial(-4163)
$ws.Range('E36').Value = '  +6.43%  '

# Row 37
$ws.Range('D37').Formula = '="0.440"'
$ws.Range('D37').Copy()
$ws.Range('D37').PasteSpecial(-4163)
$ws.Range('E37').Value = '  -1.03%  '

# Row 38
$ws.Range('D38').Value = '0.0₃0859'
$ws.Range('E38').Value = '  +4.14%  '

# Row 39
$ws.Range('D39').Formula = '="3.38"'
$ws.Range('D39').Copy()
$ws.Range('D39').PasteSpecial(-4163)
$ws.Range('E39').Value = '  -0.61%  '

# Row 40
$ws.Range('D40').Formula = '="0.147"'
$ws.Range('D40').Copy()
$ws.Range('D40').PasteSpecial(-4163)
$ws.Range('E40').Value = '  -2.78%  '

# Row 41
$ws.Range('E41').Value = '  +0.20%  '

# Row 42
$ws.Range('D42').Formula = '="11.09"'
$ws.Range('D42').Copy()
$ws.Range('D42').PasteSpecial(-4163)
$ws.Range('E42').Value = '  +12.85%  '

# Row 43
$ws.Range('D43').Formula = '="3.32"'
$ws.Range('D43').Copy()
$ws.Range('D43').PasteSpecial(-4163)
$ws.Range('E43').Value = '  -2.59%  '

# Row 44
$ws.Range('E44').Value = '  -0.14%  '

# Row 45
$ws.Range('D45').Formula = '="0.0487"'
$ws.Range('D45').Copy()
$ws.Range('D45').PasteSpecial(-4163)
$ws.Range('E45').Value = '  -0.52%  '

# Row 46
$ws.Range('D46').Formula = '="2.83"'
$ws.Range('D46').Copy()
$ws.Range('D46').PasteSpecial(-4163)
$ws.Range('E46').Value = '  +7.03%  '

# Row 47
$ws.Range('D47').Formula = '="0.149"'
$ws.Range('D47').Copy()
$ws.Range('D47').PasteSpecial(-4163)
$ws.Range('E47').Value = '  -0.67%  '

# Row 48
$ws.Range('D48').Formula = '="3.33"'
$ws.Range('D48').Copy()
$ws.Range('D48').PasteSpecial(-4163)
$ws.Range('E48').Value = '  -2.95%  '

# Row 49
$ws.Range('D49').Formula = '="3.48"'
$ws.Range('D49').Copy()
$ws.Range('D49').PasteSpecial(-4163)
$ws.Range('E49').Value = '  +7.23%  '

# Row 50
$ws.Range('D50').Formula = '="3.00"'
$ws.Range('D50').Copy()
$ws.Range('D50').PasteSpecial(-4163)
$ws.Range('E50').Value = '  -0.64%  '

# Row 51
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Formula = '="2.16"'
$ws.Range('D51').Copy()
$ws.Range('D51').PasteSpecial(-4163)
$ws.Range('E51').Value = '  +7.51%  '
